$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$txtRow3 = @'
2.3.  «К 2030 году удвоить продуктивность сельского хозяйства и доходы мелких производителей продовольствия, в частности женщин, представителей коренных народов, фермерских семейных хозяйств, скотоводов и рыбаков, (рыбоводство или рыбных хозяйств) в том числе посредством обеспечения гарантированного и равного доступа к земле, другим производственным ресурсам и факторам сельскохозяйственного производства, знаниям, финансовым услугам, рынкам и возможностям для увеличения добавленной стоимости и занятости в несельскохозяйственных секторах»
'@
$txtRow4 = @'
2.3.2.a «Доля произведенной продукции сельскохозяйственными кооперативами в общем объеме произведенной продукции»
'@
$txtRow12 = @'
Определении доли произведенной продукции сельскохозяйственными кооперативами в общем объеме произведенной продукции, в %
'@
$txtRow13 = @'
    Национальной стратегией устойчивого развития Кыргызской Республики, предусмотрены следующие меры в целях развития кооперации: использование механизма государственной помощи, при котором она будет увязана с показателями эффективности деятельности сельского товаропроизводителя - получателя помощи; при государственной поддержке аграрного сектора предпочтение будет отдаваться сельхозкооперативам, которые гарантируют сохранение уровня, а в перспективе - рост объемов производства сельскохозяйственной продукции, создание рабочих мест. При этом все формы кооперации будут иметь преимущественное право на государственную поддержку.
  Развитие сельхозкооперативов в Кыргызской Республике находится на начальном этапе и играет важную роль в укреплении экономики, конкурентоспособности, улучшении условий хозяйствования и создании стимулов для роста товарной продукции.         
Разработанная методика расчета индикатора ЦУР 2.3.2.a  будет служить информационной базой для реализации социально-экономической политики в сфере обеспечения продовольственной безопасности страны, позволит обеспечить продовольственную безопасность Кыргызстана, как важнейшую составную часть национальной    стратегией устойчивого развития, создавать условия для динамичного развития отечественного агропромышленного сектора, улучшения благосостояния населения.
'@
$txtRow14 = @'
Сельскохозяйственный кооператив который производит продукцию животноводства ежегодно один раз в год отчитываются по  статистической отчетности Форма №24(годовая) «Отчет о состоянии животноводства на 1 января»;
Сельскохозяйственный  кооператив который производит продукцию растениеводства ежегодно один раз в год отчитываются по  статистической отчетности Форма №29(годовая) «О сборе урожая сельскохозяйственных культур со всех земель и с орошаемых земель».
'@
$txtRow16 = @'
Статистическая отчетность:
- форма №24(годовая) «Отчет  о состоянии животноводства на 1 января»;
- форма №29(годовая) «О сборе урожая сельскохозяйственных культур со всех земель и с орошаемых земель».
'@
$txtRow17 = @'
 - сельскохозяйственные кооперативы (юридические лица) самостоятельно представляют отчеты по государственным формам статистической отчетности сельского хозяйства в статистические органы по территории.
  - сбор информации по отрасли растениеводства:
- сбор информации по отрасли животноводства (поголовье скота и их продукция):
'@
$txtRow19 = @'
Расчет дополнительного национального индикатора ЦУР 2.3.2.a  Доля произведенной продукции сельскохозяйственными кооперативами в общем объеме произведенной продукции в республики рассчитывается по формуле:
Дсхк=СПсхк/СВПсх*100
 , где 
СПсхк = стоимость произведенной продукции сельскохозяйственных кооперативов;
СВПсх = валовая стоимость продукции сельского хозяйства
'@

## ---- Row 2: keep text, only row height + formatting changes ----
$ws.Rows(2).RowHeight = 98.25
$ws.Range("B2").Font.Name = "Times New Roman"
$ws.Range("B2").Font.Size = 11
$ws.Range("B2").Font.Bold = $false
$ws.Range("B2").Font.Italic = $false
$ws.Range("B2").VerticalAlignment = -4108
$ws.Range("B2").WrapText = $true

## ---- Row 3: new goal text (2.3.) ----
$ws.Range("B3").Value = $txtRow3
$ws.Range("B3").Font.Name = "Times New Roman"
$ws.Range("B3").Font.Size = 11
$ws.Range("B3").Font.Bold = $false
$ws.Range("B3").Font.Italic = $false
$ws.Range("B3").VerticalAlignment = -4108
$ws.Range("B3").WrapText = $true

## ---- Row 4: new indicator text (2.3.2.a), loses left border ----
$ws.Range("B4").Value = $txtRow4
$ws.Range("B4").Font.Name = "Times New Roman"
$ws.Range("B4").Font.Size = 11
$ws.Range("B4").Font.Bold = $false
$ws.Range("B4").Font.Italic = $false
$ws.Range("B4").VerticalAlignment = -4108
$ws.Range("B4").WrapText = $true
$ws.Range("B4").Borders.Item(7).LineStyle = -4142

## ---- Row 6: text unchanged, font un-bolds / switches to Calibri ----
$ws.Range("B6").Font.Name = "Calibri"
$ws.Range("B6").Font.Size = 11
$ws.Range("B6").Font.Bold = $false
$ws.Range("B6").Font.Italic = $false
$ws.Range("B6").VerticalAlignment = -4160
$ws.Range("B6").WrapText = $true

## ---- Row 7: text unchanged, font normal Calibri, no wrap ----
$ws.Range("B7").Font.Name = "Calibri"
$ws.Range("B7").Font.Size = 11
$ws.Range("B7").Font.Bold = $false
$ws.Range("B7").Font.Italic = $false
$ws.Range("B7").VerticalAlignment = -4160
$ws.Range("B7").WrapText = $false

## ---- Row 12: new definition text ----
$ws.Range("B12").Value = $txtRow12

## ---- Row 13: new long strategy text ----
$ws.Range("B13").Value = $txtRow13

## ---- Row 14: new text + row height change ----
$ws.Rows(14).RowHeight = 76.5
$ws.Range("B14").Value = $txtRow14

## ---- Row 16: new source text ----
$ws.Range("B16").Value = $txtRow16

## ---- Row 17: new collection-method text + row height, alignment from justify to normal top ----
$ws.Rows(17).RowHeight = 148.5
$ws.Range("B17").Value = $txtRow17
$ws.Range("B17").Font.Name = "Calibri"
$ws.Range("B17").Font.Size = 11
$ws.Range("B17").Font.Bold = $false
$ws.Range("B17").Font.Italic = $false
$ws.Range("B17").HorizontalAlignment = 1
$ws.Range("B17").VerticalAlignment = -4160
$ws.Range("B17").WrapText = $true

## ---- Row 19: new calculation-method text + row height ----
$ws.Rows(19).RowHeight = 408.75
$ws.Range("B19").Value = $txtRow19

## ---- sheet view: scroll + selection to match the final editing position ----
$ws.Range("B19").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
